$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New error message strings used below
# ---------------------------------------------------------------------------
$errDuplicate = "[[400] during [POST] to [http://localhost:8090/click-escuela/school-admin/school/1234/student] [StudentController#createStudent(String,StudentApi)]: [Ya existe el estudiante]]"
$errValidation = "[[400] during [POST] to [http://localhost:8090/click-escuela/school-admin/school/1234/student] [StudentController#createStudent(String,StudentApi)]: [{`"status`":400,`"message`":`"validation error`",`"validationErrors`":[{`"message`":`"Locality cannot be null`",`"field`":`"parentApi.adressApi.locality`"},{`"message`":`"CellPhone cannot be null`",`"field`":`"parentApi.cel... (1882 bytes)]]"

# ---------------------------------------------------------------------------
# 1) New "Error" column header (U1) + highlight the existing header row green
# ---------------------------------------------------------------------------
$ws.Range("A1:T1").Interior.Color = 5287936   # FF00B050 green

$ws.Range("U1").Value = "Error"
$ws.Range("U1").Interior.Color = 1583078      # FFE62718 red
$ws.Range("U1").WrapText = $true

# ---------------------------------------------------------------------------
# 2) Error messages for the two existing rows (row 2 and row 3)
# ---------------------------------------------------------------------------
$ws.Range("U2").WrapText = $true
$ws.Range("U2").Value = $errDuplicate
$ws.Rows.Item(2).RowHeight = 165.75

$ws.Range("U3").WrapText = $true
$ws.Range("U3").Value = $errDuplicate
$ws.Rows.Item(3).RowHeight = 165.75

# ---------------------------------------------------------------------------
# 3) New student row (row 4) that failed validation
# ---------------------------------------------------------------------------
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "Mateo"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "Perez"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "98632585"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "Masculino"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Font.Underline = $true
$ws.Range("E4").Value = "2020-02-11"

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "1156968963"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Font.Underline = $true
$ws.Range("G4").Value = "tony@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:tony@gmail.com")

$ws.Range("U4").WrapText = $true
$ws.Range("U4").Value = $errValidation

$ws.Rows.Item(4).RowHeight = 331.5

# ---------------------------------------------------------------------------
# 4) Stray formatting artifacts carried over from the author's session
# ---------------------------------------------------------------------------
$ws.Range("R12").Font.Underline = $false
$ws.Range("G8").Font.Underline = $true
$ws.Range("R11").Font.Underline = $true

# ---------------------------------------------------------------------------
# 5) View state: scroll position + selection
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 16
$ws.Range("U2").Select()

Write-Output "done"
